# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Status text changes from "Ready for handoff" to "Handed back: in sync with en-US"
#    (all cells sharing that string update automatically)
#  - Per-language sheets (zh-cn, de-de) gain "Latest Target File" / "Latest Handback File"
#    links + a real "Latest Handback DateTime" for the two rows
#  - de-de handback timestamp is newer (handback ran after zh-cn)
#  - Columns that now hold longer hyperlink text are widened to fit

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1b3a336245c4fbcdfeca2c6fc701ddf30d2fa956/e2e/"

# ---- Overview sheet: just the status text + widen the status columns ----
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("E2").Value = $newStatus
$ov.Range("F2").Value = $newStatus
$ov.Range("E3").Value = $newStatus
$ov.Range("F3").Value = $newStatus
$ov.Columns.Item(5).ColumnWidth = 29.9777047293527
$ov.Columns.Item(6).ColumnWidth = 29.9777047293527

function Update-LangSheet {
    param(
        [string]$sheetName,
        [string]$handbackDateTime
    )

    $ws = $wb.Worksheets.Item($sheetName)

    # Status column (C)
    $ws.Range("C2").Value = $newStatus
    $ws.Range("C3").Value = $newStatus

    # Row 2 = 4d02d2fd... file, Row 3 = 80d98a70... file
    $md1 = "4d02d2fd-d1c2-48c9-86cc-cc0e7be7b940.md"
    $md2 = "80d98a70-eaba-4174-b056-937320179569.md"
    $xlf1 = "4d02d2fd-d1c2-48c9-86cc-cc0e7be7b940.74776b4a7d3a326d248716fa2b795f37e15e439a." + $sheetName + ".xlf"
    $xlf2 = "80d98a70-eaba-4174-b056-937320179569.7820af29c685a148564bb31da851b86f8ac2c10e." + $sheetName + ".xlf"

    # Latest Target File (I) -> hyperlink to the source .md, same target as column A
    $ws.Hyperlinks.Add($ws.Range("I2"), ($repoBase + $md1), "", "", $md1)
    $ws.Hyperlinks.Add($ws.Range("I3"), ($repoBase + $md2), "", "", $md2)

    # Latest Handback File (J) -> the generated xlf file name for this language
    $ws.Range("J2").Value = $xlf1
    $ws.Range("J3").Value = $xlf2

    # Latest Handback DateTime (K) -> now populated with the real handback time
    $ws.Range("K2").Value = $handbackDateTime
    $ws.Range("K3").Value = $handbackDateTime

    # Widen the columns that now contain longer path/hyperlink text
    $ws.Columns.Item(3).ColumnWidth = 29.9777047293527
    $ws.Columns.Item(9).ColumnWidth = 40
    $ws.Columns.Item(10).ColumnWidth = 40
}

Update-LangSheet "zh-cn" "2016-08-27 00:30:56"
Update-LangSheet "de-de" "2016-08-27 00:31:08"
